$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 292, shifting ADL..PPT (and PPT/Tahiti) down by one row.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row with the Guadalajara (GDL) colo entry.
$ws.Cells.Item(292, 1).Value = "GDL"
$ws.Cells.Item(292, 2).Value = "Guadalajara"
$ws.Cells.Item(292, 3).Value = 20.5217990875
$ws.Cells.Item(292, 4).Value = -103.3109970093
$ws.Cells.Item(292, 5).Value = "MX"
$ws.Cells.Item(292, 6).Value = "North America"
$ws.Cells.Item(292, 7).Value = "Guadalajara"

# Match the formatting used by column A on the other rows (bold, bordered, centered).
$src = $ws.Cells.Item(293, 1)
$dst = $ws.Cells.Item(292, 1)
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.Item(1).LineStyle = $src.Borders.Item(1).LineStyle
$dst.Borders.Item(2).LineStyle = $src.Borders.Item(2).LineStyle
$dst.Borders.Item(3).LineStyle = $src.Borders.Item(3).LineStyle
$dst.Borders.Item(4).LineStyle = $src.Borders.Item(4).LineStyle
